$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Exigences")
$ws.Activate()

# Rename the "Chapitre" column header (cell B1) for the export template
$ws.Range("B1").Value = "ChapitreFFFFFFFFF"

# Move the selection to the (now renamed) header cell, matching the
# author's on-screen selection when the change was made
$ws.Range("B1").Select()
